# Update column G ("K" - strikeouts) on Sheet1 with regenerated values,
# replacing the previous "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 2
    27 = 2
    28 = 1
    29 = 2
    30 = 2
    31 = 1
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    49 = 1
    50 = 0
    51 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
